$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STAGE")

# --- Clear existing content & formatting ---
$ws.Cells.Clear()

# --- Column widths (closest achievable values; column E left default) ---
$ws.Columns.Item(1).ColumnWidth = 41.59
$ws.Columns.Item(2).ColumnWidth = 18.59
$ws.Columns.Item(3).ColumnWidth = 13.6
$ws.Columns.Item(4).ColumnWidth = 15.6
$ws.Columns.Item(6).ColumnWidth = 11.6
$ws.Columns.Item(7).ColumnWidth = 14.27
$ws.Columns.Item(8).ColumnWidth = 15.59
$ws.Columns.Item(9).ColumnWidth = 34.59

# --- Row 1 header values ---
$ws.Range("A1").Value = "Location"
$ws.Range("B1").Value = "TestName"
$ws.Range("C1").Value = "User Role"
$ws.Range("D1").Value = "FirstName"
$ws.Range("E1").Value = "LastName"
$ws.Range("F1").Value = "Tenant ID"
$ws.Range("G1").Value = "password"
$ws.Range("H1").Value = "ExamName"
$ws.Range("I1").Value = "ScheduleName"

# --- Row 2 ---
$ws.Range("A2").Value = "HYD12"
$ws.Range("D2").Value = "stage"
$ws.Range("E2").Value = "controller11"
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = "@Abcd1234"
$ws.Range("H2").Value = "FPK12Exam80339"
$ws.Range("I2").Value = "FPK12Schedule15128"

# --- Row 3 ---
$ws.Range("C3").Value = "Controller"
$ws.Range("D3").Value = "fpkcontroller"
$ws.Range("E3").Value = "430836"
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = "@Abcd1234"

# --- Row 4 ---
$ws.Range("C4").Value = "Proctor"
$ws.Range("D4").Value = "FPKproctor"
$ws.Range("E4").Value = "88412"
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = "@Abcd1234"

# --- Row 5 ---
$ws.Range("C5").Value = "ExamTaker"
$ws.Range("D5").Value = "FPKexamtaker"
$ws.Range("E5").Value = "602415"
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = "@Abcd1234"

Write-Host "values written"
